$d = $word.ActiveDocument

# Find the "Date" styled paragraph (e.g. "2023-06-04") that sits right
# after the title/author block, so the Table of Contents can be inserted
# immediately after it (and before the first heading / bookmark).
$dateParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Date") {
        $dateParagraph = $p
        break
    }
}

# Collapse to just before the paragraph's end-of-paragraph mark so the
# inserted XML becomes a brand new paragraph directly after the date,
# without disturbing the date paragraph or whatever follows it.
$insertAt = $dateParagraph.Range.End - 1
$r = $d.Range($insertAt, $insertAt)

$tocXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
    + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
    + '<pkg:xmlData>' `
    + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
    + '<w:body>' `
    + '<w:sdt>' `
    + '<w:sdtPr>' `
    + '<w:docPartObj>' `
    + '<w:docPartGallery w:val="Table of Contents"/>' `
    + '<w:docPartUnique/>' `
    + '</w:docPartObj>' `
    + '</w:sdtPr>' `
    + '<w:sdtContent>' `
    + '<w:p>' `
    + '<w:pPr><w:pStyle w:val="TOCHeading"/></w:pPr>' `
    + '<w:r><w:t xml:space="preserve">Table</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">of</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">contents</w:t></w:r>' `
    + '</w:p>' `
    + '<w:p>' `
    + '<w:r>' `
    + '<w:fldChar w:fldCharType="begin" w:dirty="true"/>' `
    + '<w:instrText xml:space="preserve">TOC \o "1-3" \h \z \u</w:instrText>' `
    + '<w:fldChar w:fldCharType="separate"/>' `
    + '<w:fldChar w:fldCharType="end"/>' `
    + '</w:r>' `
    + '</w:p>' `
    + '</w:sdtContent>' `
    + '</w:sdt>' `
    + '</w:body>' `
    + '</w:document>' `
    + '</pkg:xmlData>' `
    + '</pkg:part>' `
    + '</pkg:package>'

$r.InsertXML($tocXml)
